$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(17.77300618440172, 17.45964905814496, 17.27077089761441, 17.19480542642261, 17.18225546863181, 17.26974218131344, 17.66430000732038, 18.46080315368326, 19.05300124764982, 19.32252010933319, 19.42448565929111, 19.40253154922309, 19.33091142030124, 19.28702631873198, 19.03538027117985, 18.8809534281902, 18.79214945713651, 18.76208862723874, 18.89739137542313, 19.3519514564508, 19.64842769947593, 19.49028427473223, 18.88995983910665, 18.243607671195)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}

$colC = @(10.21090539943208, 9.935730509800198, 9.767037600487511, 9.698489044850357, 9.687122239305198, 9.766112167379488, 10.11604547394275, 10.7990030673151, 11.29205890546241, 11.51319792532432, 11.59638889517166, 11.5784981141007, 11.52005363112675, 11.48418041956843, 11.27753508159243, 11.14988600864322, 11.07617393007746, 11.05116924481091, 11.16350538237288, 11.53723583168221, 11.77824822081053, 11.64994159494557, 11.15734906821109, 10.61531658243931)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}

$colD = @(14.32210881126316, 14.30993514697994, 14.30536726644407, 14.30423717562784, 14.30409370059032, 14.30534906418524, 14.31730840452092, 14.36378954189263, 14.41190648218751, 14.43680426214941, 14.44666217844202, 14.44452004360835, 14.4376066867607, 14.4334279132058, 14.41033965579413, 14.39694445367375, 14.38952333446425, 14.38705943835764, 14.39834108072125, 14.43962567275368, 14.4691102337457, 14.45314589672225, 14.39770879339201, 14.34875531475586)
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $colD[$i]
}

$colE = @(15.10668373050767, 15.11845370845852, 15.12816894143699, 15.1327528477662, 15.13355171433033, 15.1282282327872, 15.11022506917242, 15.0946997335977, 15.0953916201939, 15.09833890839526, 15.09983361792178, 15.09949486623686, 15.09845429081955, 15.09786621710036, 15.09525199235641, 15.09432262121267, 15.09403586591347, 15.0939813290593, 15.0943959074428, 15.09874965725395, 15.103801970521, 15.10090354626162, 15.09436200367526, 15.09677742828761)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $colE[$i]
}

$colG = @(61.3032098326869, 60.60837399812401, 60.1932914575972, 60.027197844362, 59.99980697753114, 60.1910388933179, 61.06132656648183, 62.85260696429496, 64.21054086203142, 64.83525535811164, 65.07264105272411, 65.02148226188721, 64.85476984490765, 64.75275498307667, 64.16983906314412, 63.81389648893955, 63.60983612190822, 63.5408652036361, 63.85171934441441, 64.90371653872116, 65.59595144236559, 65.22612387820772, 63.8346178241035, 62.3599436967169)
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $colG[$i]
}

$colH = @(21.81531885108973, 21.75340912045176, 21.71953317597767, 21.70677177876887, 21.7047158416569, 21.71935684383637, 21.79311591484332, 21.97040401039854, 22.12026163968668, 22.19260791223862, 22.22059505293245, 22.21454138561998, 22.19489863930894, 22.18294358865947, 22.11561689159868, 22.07537727946336, 22.05262613134956, 22.04499085194504, 22.07962018124222, 22.20065222389637, 22.28319455178505, 22.23882860776582, 22.07770077126133, 21.91897445418082)
for ($i = 0; $i -lt $colH.Length; $i++) {
    $ws.Cells.Item($i + 2, 8).Value = $colH[$i]
}

$colJ = @(8.925410521891838, 8.942316113811593, 8.95337378833463, 8.958050628361086, 8.958837537646721, 8.953436170003167, 8.931099160207616, 8.892656102707745, 8.867656722053837, 8.856983643494139, 8.853042218859981, 8.853886622005989, 8.856657372976311, 8.858367582696305, 8.86836827682585, 8.874682251919481, 8.878379720116717, 8.879642936491235, 8.874003306937798, 8.855840818294286, 8.844554706035828, 8.850524972102759, 8.874310047755612, 8.902484533293695)
for ($i = 0; $i -lt $colJ.Length; $i++) {
    $ws.Cells.Item($i + 2, 10).Value = $colJ[$i]
}

$colM = @(20.41769206574797, 20.34929384773791, 20.31254709294975, 20.29890314649385, 20.29671822329945, 20.31235768554026, 20.39302548359972, 20.5923397006558, 20.76303546693301, 20.84575419257045, 20.87778721868203, 20.87085711449946, 20.8483754712034, 20.83469657734598, 20.75772994630353, 20.71179800409764, 20.68585747880272, 20.67715718962504, 20.71663817213086, 20.85495978836286, 20.94948503038105, 20.89866455741179, 20.71444847808503, 20.5340950783101)
for ($i = 0; $i -lt $colM.Length; $i++) {
    $ws.Cells.Item($i + 2, 13).Value = $colM[$i]
}
